# Rename the "_old" / "_new" header-suffix column headers to the
# respective AHB format-version names ("_FV2210" / "_FV2304"), then turn
# the header row + data into an Excel Table (ListObject) and freeze the
# header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header row (row 1) cells ----------------------------------
# Columns A-J carried the "_old" suffix, columns L-U carried the "_new"
# suffix (column K is the unchanged "diff" column).
$oldSuffixCols = @("A","B","C","D","E","F","G","H","I","J")
$newSuffixCols = @("L","M","N","O","P","Q","R","S","T","U")

$oldBaseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")
$newBaseNames = $oldBaseNames

for ($i = 0; $i -lt $oldSuffixCols.Length; $i++) {
    $ws.Range($oldSuffixCols[$i] + "1").Value = $oldBaseNames[$i] + "_FV2210"
}

for ($i = 0; $i -lt $newSuffixCols.Length; $i++) {
    $ws.Range($newSuffixCols[$i] + "1").Value = $newBaseNames[$i] + "_FV2304"
}

# --- 2) Convert the data range into an Excel Table (ListObject) ----------
$dataRange = $ws.Range("A1:U89")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"

# --- 3) Freeze the header row ---------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Header rename + table + freeze panes applied"
